$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.320.13"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "3.433.90"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'413.79"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "'129.89"
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.724"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").Value = "'42.66"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "'9.42"
$ws.Range("E12").Value = "  +3.20%  "
$ws.Range("E13").Value = "  +5.71%  "
$ws.Range("D14").Value = "3.977.70"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  -3.47%  "
$ws.Range("D17").Value = "'12.89"
$ws.Range("E17").Value = "  +4.00%  "
$ws.Range("D18").Value = "3.405.27"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "62.372.00"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "'476.37"
$ws.Range("E21").Value = "  +7.33%  "
$ws.Range("D22").Value = "'91.25"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'3.29"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "'13.41"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").Value = "'10.59"
$ws.Range("E25").Value = "  +21.97%  "
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("D27").Value = "'33.24"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'11.93"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "'2.66"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("D34").Value = "'40.64"
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'58.44"
$ws.Range("E36").Value = "  +9.99%  "
$ws.Range("D37").Value = "'0.0489"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'3.03"
$ws.Range("E39").Value = "  +3.26%  "
$ws.Range("E40").Value = "  +4.36%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").Value = "'145.42"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").Value = "'2.67"
$ws.Range("E44").Value = "  +6.40%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'4.36"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'2.06"
$ws.Range("E46").Value = "  +3.88%  "
$ws.Range("E47").Value = "  +13.46%  "
$ws.Range("B48").Value = "Celestia"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D48").Value = "'16.37"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.0₃0543"
$ws.Range("E49").Value = "  +32.89%  "
$ws.Range("D50").Value = "'22.38"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.140"
$ws.Range("E51").Value = "  +4.55%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
